$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the SNAP2 chassis # value for the first data row (LWA-213)
$ws.Range("K4").Value = 3

# Slightly narrow the columns (pretest width tweak)
$ws.Range("A1").EntireColumn.ColumnWidth = 5.646258503401357
$ws.Range("B1").EntireColumn.ColumnWidth = 6.860544217687077
$ws.Range("C1").EntireColumn.ColumnWidth = 9.697278911564666
$ws.Range("D1").EntireColumn.ColumnWidth = 11.314625850340166
$ws.Range("E1").EntireColumn.ColumnWidth = 9.559523809523766
$ws.Range("F1:G1").EntireColumn.ColumnWidth = 6.727891156462587
$ws.Range("H1").EntireColumn.ColumnWidth = 8.482993197278907
$ws.Range("I1").EntireColumn.ColumnWidth = 17.253401360544267
$ws.Range("J1").EntireColumn.ColumnWidth = 5.778911564625847
$ws.Range("K1").EntireColumn.ColumnWidth = 7.130952380952377
$ws.Range("L1:M1").EntireColumn.ColumnWidth = 7.671768707482998
$ws.Range("N1").EntireColumn.ColumnWidth = 8.753401360544217
$ws.Range("O1").EntireColumn.ColumnWidth = 5.243197278911567
$ws.Range("P1").EntireColumn.ColumnWidth = 13.748299319727867
$ws.Range("Q1").EntireColumn.ColumnWidth = 8.886054421768707
$ws.Range("R1:S1").EntireColumn.ColumnWidth = 8.212585034013607
$ws.Range("T1").EntireColumn.ColumnWidth = 6.590136054421767
$ws.Range("U1:V1").EntireColumn.ColumnWidth = 6.998299319727887
$ws.Range("W1").EntireColumn.ColumnWidth = 8.753401360544217
$ws.Range("X1").EntireColumn.ColumnWidth = 8.482993197278907
$ws.Range("Y1:Z1").EntireColumn.ColumnWidth = 8.753401360544217
$ws.Range("AA1").EntireColumn.ColumnWidth = 28.595238095238066
